$d = $word.ActiveDocument
$pp = $d.Paragraphs(39)
$insertPos = $pp.Range.End - 1
$ins = $d.Range($insertPos, $insertPos)
$ins.InsertAfter("  -- object consist of (menuItem, price) you have to take for customised menu.")
Write-Output "done"
$p2 = $d.Paragraphs(39)
Write-Output $p2.Range.Text
